$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.604.19'
$ws.Cells.Item(2, 5).Value = '  +1.07%  '
$ws.Cells.Item(3, 4).Value = '1.881.94'
$ws.Cells.Item(3, 5).Value = '  +0.04%  '
$ws.Cells.Item(4, 4).Value = '0.9986'
$ws.Cells.Item(4, 5).Value = '  -0.23%  '
$ws.Cells.Item(5, 4).Value = '239.31'
$ws.Cells.Item(5, 5).Value = '  +0.88%  '
$ws.Cells.Item(6, 4).Value = '0.9993'
$ws.Cells.Item(6, 5).Value = '  -0.16%  '
$ws.Cells.Item(7, 4).Value = '0.4810'
$ws.Cells.Item(7, 5).Value = '  +0.05%  '
$ws.Cells.Item(8, 4).Value = '0.2836'
$ws.Cells.Item(8, 5).Value = '  -1.52%  '
$ws.Cells.Item(9, 4).Value = '0.06547'
$ws.Cells.Item(9, 5).Value = '  -0.69%  '
$ws.Cells.Item(10, 4).Value = '2.324.35'
$ws.Cells.Item(10, 5).Value = '  +22.69%  '
$ws.Cells.Item(11, 4).Value = '0.07486'
$ws.Cells.Item(12, 4).Value = '16.70'
$ws.Cells.Item(12, 5).Value = '  -0.94%  '
$ws.Cells.Item(13, 4).Value = '5.105'
$ws.Cells.Item(13, 5).Value = '  -1.12%  '
$ws.Cells.Item(14, 4).Value = '88.62'
$ws.Cells.Item(14, 5).Value = '  +1.20%  '
$ws.Cells.Item(15, 4).Value = '0.6665'
$ws.Cells.Item(15, 5).Value = '  +0.93%  '
$ws.Cells.Item(16, 4).Value = '30.550.98'
$ws.Cells.Item(16, 5).Value = '  +0.98%  '
$ws.Cells.Item(17, 4).Value = '13.37'
$ws.Cells.Item(17, 5).Value = '  -0.42%  '
$ws.Cells.Item(18, 4).Value = '0.9989'
$ws.Cells.Item(18, 5).Value = '  -0.19%  '
$ws.Cells.Item(19, 4).Value = '0.000007627'
$ws.Cells.Item(19, 5).Value = '  -1.33%  '
$ws.Cells.Item(20, 4).Value = '2.124.84'
$ws.Cells.Item(20, 5).Value = '  -1.10%  '
$ws.Cells.Item(21, 4).Value = '5.347'
$ws.Cells.Item(21, 5).Value = '  -1.85%  '
$ws.Cells.Item(22, 4).Value = '224.47'
$ws.Cells.Item(22, 5).Value = '  +16.55%  '
$ws.Cells.Item(23, 5).Value = '  -0.22%  '
$ws.Cells.Item(24, 2).Value = 'BitDAO'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Cells.Item(24, 4).Value = '0.3987'
$ws.Cells.Item(24, 5).Value = '  -8.16%  '
$ws.Cells.Item(25, 2).Value = 'Chainlink'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(25, 4).Value = '6.225'
$ws.Cells.Item(25, 5).Value = '  +0.85%  '
$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).Value = '9.371'
$ws.Cells.Item(26, 5).Value = '  -0.43%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).Value = '166.74'
$ws.Cells.Item(27, 5).Value = '  -0.18%  '
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).Value = '18.63'
$ws.Cells.Item(28, 5).Value = '  +1.74%  '
$ws.Cells.Item(29, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(29, 4).Value = '1.979'
$ws.Cells.Item(29, 5).Value = '  +2.13%  '
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(30, 4).Value = '1.455'
$ws.Cells.Item(30, 5).Value = '  +0.70%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).Value = '0.09471'
$ws.Cells.Item(31, 5).Value = '  +3.45%  '
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '4.339'
$ws.Cells.Item(32, 5).Value = '  +1.85%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).Value = '4.039'
$ws.Cells.Item(33, 5).Value = '  -0.17%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '0.05040'
$ws.Cells.Item(34, 5).Value = '  -0.41%  '
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).Value = '1.210'
$ws.Cells.Item(35, 5).Value = '  +5.88%  '
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(36, 4).Value = '0.7501'
$ws.Cells.Item(36, 5).Value = '  +1.81%  '
$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(37, 4).Value = '2.699'
$ws.Cells.Item(37, 5).Value = '  -0.49%  '
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '0.01839'
$ws.Cells.Item(38, 5).Value = '  +0.69%  '
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(39, 4).Value = '2.621'
$ws.Cells.Item(39, 5).Value = '  -0.60%  '
$ws.Cells.Item(40, 4).Value = '2.095'
$ws.Cells.Item(40, 5).Value = '  +1.23%  '
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).Value = '0.9085'
$ws.Cells.Item(41, 5).Value = '  -0.47%  '
$ws.Cells.Item(42, 2).Value = 'Quant'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(42, 4).Value = '106.19'
$ws.Cells.Item(42, 5).Value = '  -0.19%  '
$ws.Cells.Item(43, 4).Value = '0.4299'
$ws.Cells.Item(43, 5).Value = '  -0.50%  '
$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).Value = '5.841'
$ws.Cells.Item(44, 5).Value = '  -0.77%  '
$ws.Cells.Item(45, 2).Value = 'PaxDollar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(45, 4).Value = '1.006'
$ws.Cells.Item(45, 5).Value = '  +0.60%  '
$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).Value = '7.502'
$ws.Cells.Item(46, 5).Value = '  -1.97%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).Value = '64.84'
$ws.Cells.Item(47, 5).Value = '  -0.24%  '
$ws.Cells.Item(48, 2).Value = 'Algorand'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(48, 4).Value = '0.1285'
$ws.Cells.Item(48, 5).Value = '  -4.97%  '
$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49, 4).Value = '1.481'
$ws.Cells.Item(49, 5).Value = '  -6.75%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '8.936'
$ws.Cells.Item(50, 5).Value = '  +0.44%  '
$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(51, 4).Value = '33.80'
$ws.Cells.Item(51, 5).Value = '  -1.15%  '
